$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 'Malik Beasley'
$ws.Range("B5").Value = 'SG'
$ws.Range("C5").Value = 'Detroit Pistons'

$ws.Range("A6").Value = 'Josh Hart'
$ws.Range("B6").Value = 'SF,PF'
$ws.Range("C6").Value = 'New York Knicks'

$ws.Range("A7").Value = 'Tari Eason'
$ws.Range("B7").Value = 'SF,PF'
$ws.Range("C7").Value = 'Houston Rockets'

$ws.Range("A8").Value = 'Alexandre Sarr'
$ws.Range("B8").Value = 'PF,C'
$ws.Range("C8").Value = 'Washington Wizards'

$ws.Range("A9").Value = 'Toumani Camara'
$ws.Range("B9").Value = 'SF,PF'
$ws.Range("C9").Value = 'Portland Trail Blazers'

$ws.Range("A10").Value = 'Grant Williams'
$ws.Range("B10").Value = 'PF,C'
$ws.Range("C10").Value = 'Charlotte Hornets'

$ws.Range("A11").Value = 'Victor Wembanyama'
$ws.Range("B11").Value = 'C'
$ws.Range("C11").Value = 'San Antonio Spurs'

$ws.Range("A12").Value = 'Myles Turner'
$ws.Range("B12").Value = 'C'
$ws.Range("C12").Value = 'Indiana Pacers'

$ws.Range("A13").Value = 'Domantas Sabonis'
$ws.Range("B13").Value = 'C'
$ws.Range("C13").Value = 'Sacramento Kings'

$ws.Range("A14").Value = 'De''Andre Hunter'
$ws.Range("B14").Value = 'SF,PF'
$ws.Range("C14").Value = 'Atlanta Hawks'

$ws.Range("A15").Value = 'Naz Reid'
$ws.Range("B15").Value = 'PF,C'
$ws.Range("C15").Value = 'Minnesota Timberwolves'

$ws.Range("A16").Value = 'Donovan Mitchell'
$ws.Range("B16").Value = 'PG,SG'
$ws.Range("C16").Value = 'Cleveland Cavaliers'

$ws.Range("A17").Value = 'Bradley Beal'
$ws.Range("B17").Value = 'PG,SG,SF'
$ws.Range("C17").Value = 'Phoenix Suns'

$ws.Range("A18").Value = 'Kristaps Porzingis'
$ws.Range("B18").Value = 'PF,C'
$ws.Range("C18").Value = 'Boston Celtics'

$ws.Range("A19").Value = 'Ja Morant'
$ws.Range("B19").Value = 'PG'
$ws.Range("C19").Value = 'Memphis Grizzlies'
